$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.249.31"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "3.182.33"
$ws.Range("E3").Value = "  -7.53%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.67"
$ws.Range("E5").Value = "  -3.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.34"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.606"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "3.177.92"
$ws.Range("E9").Value = "  -7.70%  "
$ws.Range("E10").Value = "  -5.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.63"
$ws.Range("E11").Value = "  -4.26%  "
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("D13").Value = "3.728.47"
$ws.Range("E13").Value = "  -7.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.136"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.55"
$ws.Range("E15").Value = "  -4.62%  "
$ws.Range("D16").Value = "64.235.17"
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000163"
$ws.Range("E17").Value = "  -4.54%  "
$ws.Range("D18").Value = "3.179.29"
$ws.Range("E18").Value = "  -7.95%  "
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("E20").Value = "  -4.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "354.10"
$ws.Range("E21").Value = "  -4.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.23"
$ws.Range("E22").Value = "  -4.63%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.97"
$ws.Range("E24").Value = "  -4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.507"
$ws.Range("E25").Value = "  -4.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.14"
$ws.Range("E33").Value = "  -5.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.68"
$ws.Range("E34").Value = "  -4.18%  "
$ws.Range("E35").Value = "  -5.15%  "
$ws.Range("E36").Value = "  -5.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "155.08"
$ws.Range("E37").Value = "  -3.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.808"
$ws.Range("E38").Value = "  -8.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.98"
$ws.Range("E39").Value = "  -8.40%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.51"
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("D42").Value = "2.656.15"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.18"
$ws.Range("E43").Value = "  -5.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.01"
$ws.Range("E44").Value = "  -7.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "332.71"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0656"
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.75"
$ws.Range("E49").Value = "  -5.95%  "
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("E51").Value = "  -0.14%  "
